# Fruta / hortaliza, semanal
#
# The weekly refresh reshuffles the per-record fields (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg
# and Kg/unidad) across the existing data rows (2-16) of the sheet. The
# "dimension" columns (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID,
# Producto, Categoria ID, Categoria, Variedad) stay put since they are already
# identical for every row.
#
# Because the new value for one row may come from another row that itself is
# being overwritten, every original value used below is captured first into
# an in-memory table, and only then are the cells written back out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose contents get redistributed across rows 2-16.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Maps each row to the row whose (D,L,M,N,O,P,Q,R,S,T) block it receives.
$perm = @{
    2  = 5
    3  = 11
    4  = 12
    5  = 9
    6  = 10
    7  = 3
    8  = 2
    9  = 13
    10 = 14
    11 = 4
    12 = 16
    13 = 15
    14 = 7
    15 = 8
    16 = 6
}

# Snapshot every current value before any cell gets overwritten.
$orig = @{}
foreach ($r in 2..16) {
    foreach ($c in $cols) {
        $orig["$c$r"] = $ws.Range("$c$r").Value2
    }
}

# Write the permuted values back into the sheet.
foreach ($r in 2..16) {
    $src = $perm[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig["$c$src"]
    }
}
